$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 400, shifting existing data (rows 400-448) down to 403-451
$ws.Rows("400:402").Insert()

# Populate new row 400
$ws.Range("A400").Value = 7
$ws.Range("B400").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C400").Value = 'Ñuble'
$ws.Range("D400").Value = 45131
$ws.Range("E400").Value = 16
$ws.Range("F400").Value = 'Fruta'
$ws.Range("G400").Value = 100104
$ws.Range("H400").Value = 'Frutos de pepita'
$ws.Range("I400").Value = 100104005
$ws.Range("J400").Value = 'Pera'
$ws.Range("K400").Value = 'Abate Fettel'
$ws.Range("L400").Value = 'Primera'
$ws.Range("M400").Value = 40
$ws.Range("N400").Value = 10000
$ws.Range("O400").Value = 10000
$ws.Range("P400").Value = 10000
$ws.Range("Q400").Value = '$/bandeja 18 kilos granel'
$ws.Range("R400").Value = 'Región de O''Higgins'
$ws.Range("S400").Value = 556
$ws.Range("T400").Value = 18

# Populate new row 401
$ws.Range("A401").Value = 7
$ws.Range("B401").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C401").Value = 'Ñuble'
$ws.Range("D401").Value = 45131
$ws.Range("E401").Value = 16
$ws.Range("F401").Value = 'Fruta'
$ws.Range("G401").Value = 100104
$ws.Range("H401").Value = 'Frutos de pepita'
$ws.Range("I401").Value = 100104005
$ws.Range("J401").Value = 'Pera'
$ws.Range("K401").Value = 'Packham''s Triumph'
$ws.Range("L401").Value = 'Especial'
$ws.Range("M401").Value = 60
$ws.Range("N401").Value = 12000
$ws.Range("O401").Value = 12000
$ws.Range("P401").Value = 12000
$ws.Range("Q401").Value = '$/bandeja 18 kilos granel'
$ws.Range("R401").Value = 'Región de O''Higgins'
$ws.Range("S401").Value = 667
$ws.Range("T401").Value = 18

# Populate new row 402
$ws.Range("A402").Value = 7
$ws.Range("B402").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C402").Value = 'Ñuble'
$ws.Range("D402").Value = 45131
$ws.Range("E402").Value = 16
$ws.Range("F402").Value = 'Fruta'
$ws.Range("G402").Value = 100104
$ws.Range("H402").Value = 'Frutos de pepita'
$ws.Range("I402").Value = 100104005
$ws.Range("J402").Value = 'Pera'
$ws.Range("K402").Value = 'Packham''s Triumph'
$ws.Range("L402").Value = 'Primera'
$ws.Range("M402").Value = 120
$ws.Range("N402").Value = 10000
$ws.Range("O402").Value = 10000
$ws.Range("P402").Value = 10000
$ws.Range("Q402").Value = '$/bandeja 18 kilos granel'
$ws.Range("R402").Value = 'Región de O''Higgins'
$ws.Range("S402").Value = 556
$ws.Range("T402").Value = 18

Write-Output "Done"